# Auto-generated script to restructure BMED prerequisites sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nbsp = [char]0x00A0

# --- Update header row (row 1): insert Corequisites, Concurrent, Recommended before "Terms Typically Offered" ---
$ws.Range("G1").Value = $ws.Range("D1").Text
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# --- Process each data row (2-58): move old column D (Terms Typically Offered) to G, ---
# --- set D/E/F defaults to "NA", then apply row-specific corrections below ---
for ($r = 2; $r -le 58; $r++) {
    $oldD = $ws.Cells.Item($r, 4).Text
    $ws.Cells.Item($r, 7).Value = $oldD
    $ws.Cells.Item($r, 4).Value = "NA"
    $ws.Cells.Item($r, 5).Value = "NA"
    $ws.Cells.Item($r, 6).Value = "NA"
    $c = $ws.Cells.Item($r, 3).Text
    $ws.Cells.Item($r, 3).Value = $c.Replace($nbsp, " ")
}

# --- Row-specific corrections: split Corequisite/Recommended notes out of Prerequisites text, ---
# --- and fix up the moved "Terms Typically Offered" text where it had concatenated notes ---

# Row 5
$ws.Cells.Item(5, 3).Value = "MATH 142; for engineering students only."
$ws.Cells.Item(5, 4).Value = "BIO 213."
$ws.Cells.Item(5, 5).Value = "NA"
$ws.Cells.Item(5, 6).Value = "CHEM 124."
$ws.Cells.Item(5, 7).Value = "F,W,SP,SU  "

# Row 8
$ws.Cells.Item(8, 3).Value = "MATE 210, ME 328, STAT 312."
$ws.Cells.Item(8, 4).Value = "NA"
$ws.Cells.Item(8, 5).Value = "NA"
$ws.Cells.Item(8, 6).Value = "BMED 420, BMED 460."
$ws.Cells.Item(8, 7).Value = "SP "

# Row 14
$ws.Cells.Item(14, 3).Value = "CE 204 or CE 208; and ME 212."
$ws.Cells.Item(14, 4).Value = "BMED 310."
$ws.Cells.Item(14, 5).Value = "NA"
$ws.Cells.Item(14, 6).Value = "NA"
$ws.Cells.Item(14, 7).Value = "W, SP "

# Row 15
$ws.Cells.Item(15, 3).Value = "CE 204 or CE 208; and MATE 210."
$ws.Cells.Item(15, 4).Value = "BMED 310."
$ws.Cells.Item(15, 5).Value = "NA"
$ws.Cells.Item(15, 6).Value = "NA"
$ws.Cells.Item(15, 7).Value = "W, SP "

# Row 21
$ws.Cells.Item(21, 3).Value = "NA"
$ws.Cells.Item(21, 4).Value = "BMED 434/EE 423/MATE 430."
$ws.Cells.Item(21, 5).Value = "NA"
$ws.Cells.Item(21, 6).Value = "NA"
$ws.Cells.Item(21, 7).Value = "W"

# Row 38
$ws.Cells.Item(38, 3).Value = "one of the ASCI 438, BIO 361, or BMED 460; or graduate standing."
$ws.Cells.Item(38, 4).Value = "NA"
$ws.Cells.Item(38, 5).Value = "NA"
$ws.Cells.Item(38, 6).Value = "NA"
$ws.Cells.Item(38, 7).Value = "F"

# Row 50
$ws.Cells.Item(50, 3).Value = "ASCI 438, BIO 361, or BMED 460; and STAT 218 or STAT 312."
$ws.Cells.Item(50, 4).Value = "BMED 560."
$ws.Cells.Item(50, 5).Value = "NA"
$ws.Cells.Item(50, 6).Value = "NA"
$ws.Cells.Item(50, 7).Value = "SP "
